$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename SKU codes from the old "KRTD" prefix to the new "PIPI" prefix.
$ws.Range("A2").Value = "PIPI00001"
$ws.Range("A3").Value = "PIPI00002"
$ws.Range("A4").Value = "PIPI00003"
$ws.Range("A5").Value = "PIPI00004"
$ws.Range("A6").Value = "PIPI00005"

# Move / extend the current selection to the SKU column (A2:A6), active cell A2.
$ws.Range("A2:A6").Select()
